$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teamUpdates = @{
    'D14' = 'Oklahoma City Thunder'
    'C22' = 'Oklahoma City Thunder'
    'C30' = 'Brooklyn Nets'
    'C40' = 'Oklahoma City Thunder'
    'C44' = 'Brooklyn Nets'
    'C51' = 'Oklahoma City Thunder'
    'D55' = 'Brooklyn Nets'
    'D66' = 'Oklahoma City Thunder'
    'D69' = 'Brooklyn Nets'
    'C75' = 'Oklahoma City Thunder'
    'C90' = 'Brooklyn Nets'
    'C93' = 'Oklahoma City Thunder'
    'D97' = 'Oklahoma City Thunder'
    'C106' = 'Brooklyn Nets'
    'C113' = 'Oklahoma City Thunder'
    'C119' = 'Brooklyn Nets'
    'D125' = 'Oklahoma City Thunder'
    'D142' = 'Brooklyn Nets'
    'C143' = 'Oklahoma City Thunder'
    'D156' = 'Brooklyn Nets'
    'C164' = 'Oklahoma City Thunder'
    'D170' = 'Brooklyn Nets'
    'D173' = 'Oklahoma City Thunder'
    'C174' = 'Brooklyn Nets'
    'D183' = 'Oklahoma City Thunder'
    'C193' = 'Brooklyn Nets'
    'C201' = 'Oklahoma City Thunder'
    'C204' = 'Brooklyn Nets'
    'D212' = 'Brooklyn Nets'
    'C219' = 'Oklahoma City Thunder'
    'D225' = 'Brooklyn Nets'
    'C230' = 'Oklahoma City Thunder'
    'D234' = 'Brooklyn Nets'
    'D237' = 'Oklahoma City Thunder'
    'C254' = 'Brooklyn Nets'
    'D254' = 'Oklahoma City Thunder'
    'C275' = 'Brooklyn Nets'
    'C281' = 'Oklahoma City Thunder'
    'C294' = 'Brooklyn Nets'
    'C295' = 'Oklahoma City Thunder'
    'C305' = 'Brooklyn Nets'
    'D310' = 'Brooklyn Nets'
    'C318' = 'Oklahoma City Thunder'
    'C330' = 'Brooklyn Nets'
    'C333' = 'Oklahoma City Thunder'
    'D343' = 'Brooklyn Nets'
    'C355' = 'Oklahoma City Thunder'
    'C359' = 'Brooklyn Nets'
    'D370' = 'Brooklyn Nets'
    'D372' = 'Oklahoma City Thunder'
    'D378' = 'Oklahoma City Thunder'
    'C400' = 'Brooklyn Nets'
    'C406' = 'Brooklyn Nets'
    'D408' = 'Oklahoma City Thunder'
    'D416' = 'Brooklyn Nets'
    'C423' = 'Oklahoma City Thunder'
    'C429' = 'Brooklyn Nets'
    'C440' = 'Brooklyn Nets'
    'D444' = 'Oklahoma City Thunder'
    'D456' = 'Brooklyn Nets'
    'C457' = 'Oklahoma City Thunder'
    'C471' = 'Oklahoma City Thunder'
    'D471' = 'Brooklyn Nets'
    'D480' = 'Brooklyn Nets'
    'C484' = 'Oklahoma City Thunder'
    'C493' = 'Brooklyn Nets'
    'D499' = 'Oklahoma City Thunder'
    'D504' = 'Oklahoma City Thunder'
    'D512' = 'Brooklyn Nets'
    'C521' = 'Oklahoma City Thunder'
    'C531' = 'Brooklyn Nets'
    'D540' = 'Oklahoma City Thunder'
    'C551' = 'Brooklyn Nets'
    'D554' = 'Oklahoma City Thunder'
    'D562' = 'Oklahoma City Thunder'
    'C566' = 'Brooklyn Nets'
    'D572' = 'Brooklyn Nets'
    'C574' = 'Oklahoma City Thunder'
    'C588' = 'Brooklyn Nets'
    'D592' = 'Oklahoma City Thunder'
    'D604' = 'Oklahoma City Thunder'
    'D609' = 'Brooklyn Nets'
    'D617' = 'Oklahoma City Thunder'
    'D621' = 'Brooklyn Nets'
    'D628' = 'Oklahoma City Thunder'
    'D636' = 'Brooklyn Nets'
    'D640' = 'Oklahoma City Thunder'
    'D646' = 'Brooklyn Nets'
    'D653' = 'Oklahoma City Thunder'
    'C662' = 'Brooklyn Nets'
    'C677' = 'Brooklyn Nets'
    'C682' = 'Oklahoma City Thunder'
    'C688' = 'Brooklyn Nets'
    'D698' = 'Oklahoma City Thunder'
    'C713' = 'Oklahoma City Thunder'
    'C717' = 'Brooklyn Nets'
    'D726' = 'Brooklyn Nets'
    'C729' = 'Oklahoma City Thunder'
    'D737' = 'Brooklyn Nets'
    'C745' = 'Oklahoma City Thunder'
    'D759' = 'Oklahoma City Thunder'
    'C760' = 'Brooklyn Nets'
    'D765' = 'Brooklyn Nets'
    'D773' = 'Oklahoma City Thunder'
    'C779' = 'Brooklyn Nets'
    'C788' = 'Oklahoma City Thunder'
    'C792' = 'Brooklyn Nets'
    'D804' = 'Brooklyn Nets'
    'D806' = 'Oklahoma City Thunder'
    'C816' = 'Brooklyn Nets'
    'C819' = 'Oklahoma City Thunder'
    'C832' = 'Brooklyn Nets'
    'C839' = 'Oklahoma City Thunder'
    'D847' = 'Brooklyn Nets'
    'C855' = 'Oklahoma City Thunder'
    'C871' = 'Brooklyn Nets'
    'D876' = 'Oklahoma City Thunder'
    'D878' = 'Brooklyn Nets'
    'D882' = 'Oklahoma City Thunder'
    'C898' = 'Oklahoma City Thunder'
    'D901' = 'Brooklyn Nets'
    'D914' = 'Oklahoma City Thunder'
    'D917' = 'Oklahoma City Thunder'
    'C919' = 'Brooklyn Nets'
    'D928' = 'Brooklyn Nets'
    'C935' = 'Oklahoma City Thunder'
    'D944' = 'Brooklyn Nets'
    'D945' = 'Oklahoma City Thunder'
    'C953' = 'Brooklyn Nets'
    'C962' = 'Oklahoma City Thunder'
    'C975' = 'Oklahoma City Thunder'
    'D991' = 'Oklahoma City Thunder'
    'C992' = 'Brooklyn Nets'
    'D997' = 'Brooklyn Nets'
    'C1006' = 'Oklahoma City Thunder'
    'D1014' = 'Brooklyn Nets'
    'D1016' = 'Oklahoma City Thunder'
    'D1025' = 'Oklahoma City Thunder'
    'D1040' = 'Brooklyn Nets'
    'C1045' = 'Oklahoma City Thunder'
    'D1048' = 'Brooklyn Nets'
    'C1067' = 'Oklahoma City Thunder'
    'D1071' = 'Brooklyn Nets'
    'D1080' = 'Oklahoma City Thunder'
    'D1084' = 'Brooklyn Nets'
    'D1091' = 'Oklahoma City Thunder'
    'D1092' = 'Brooklyn Nets'
    'D1114' = 'Brooklyn Nets'
    'C1124' = 'Brooklyn Nets'
    'C1126' = 'Oklahoma City Thunder'
    'D1133' = 'Oklahoma City Thunder'
    'C1141' = 'Brooklyn Nets'
    'C1146' = 'Oklahoma City Thunder'
    'C1157' = 'Brooklyn Nets'
    'D1162' = 'Oklahoma City Thunder'
    'D1169' = 'Brooklyn Nets'
    'D1176' = 'Oklahoma City Thunder'
    'D1177' = 'Brooklyn Nets'
    'D1189' = 'Oklahoma City Thunder'
    'D1198' = 'Brooklyn Nets'
    'C1206' = 'Brooklyn Nets'
    'C1210' = 'Oklahoma City Thunder'
    'C1220' = 'Oklahoma City Thunder'
    'C1222' = 'Brooklyn Nets'
    'C1233' = 'Brooklyn Nets'
    'C1238' = 'Oklahoma City Thunder'
    'C1239' = 'Brooklyn Nets'
    'C1244' = 'Oklahoma City Thunder'
    'D1248' = 'Brooklyn Nets'
    'D1253' = 'Brooklyn Nets'
    'D1256' = 'Oklahoma City Thunder'
    'C1262' = 'Brooklyn Nets'
    'D1263' = 'Oklahoma City Thunder'
    'C1268' = 'Oklahoma City Thunder'
    'D1269' = 'Brooklyn Nets'
    'D1273' = 'Oklahoma City Thunder'
    'C1275' = 'Brooklyn Nets'
    'C1277' = 'Oklahoma City Thunder'
    'C1281' = 'Oklahoma City Thunder'
    'D1286' = 'Oklahoma City Thunder'
    'D1290' = 'Oklahoma City Thunder'
    'C1294' = 'Oklahoma City Thunder'
}

foreach ($cell in $teamUpdates.Keys) {
    $ws.Range($cell).Value = $teamUpdates[$cell]
}
